$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biology")
$ws.Range("B4").Value = "octopus"
$ws.Range("B2").Value = "skin"
$ws.Range("B11").Select()
